$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "undo/redo" column (C) previously marked row 7 ("5. set current player by
# playerID") as implemented ("Y"). Clear it to indicate undo/redo was NOT
# implemented for this item.
$ws.Range("C7").ClearContents()

# Column C widened (the comment/marker column grew a bit after the edit).
$ws.Columns("C").ColumnWidth = 21.857142857142854

# Move the active selection to E6, matching where the author left off editing.
$ws.Range("E6").Select()
